# Update automatico via Actualizar 02-05-2021 22-01-25
# Appends one more "pull cycle" of 14 rows (828-841) to the availability
# log sheet, following the same 14-row repeating pattern already present
# in the sheet, and nudges the timestamp recorded for the previous cycle
# (rows 814-827) by a few micro-days (artifact of two check runs that
# landed inside the same polling pass).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) The previous cycle (rows 814-827) got its "last checked" timestamp
#    refined very slightly (44232.89639873675 -> 44232.89639873843).
# ---------------------------------------------------------------------
for ($r = 814; $r -le 827; $r++) {
    $ws.Cells.Item($r, 4).Value = 44232.89639873843
}

# ---------------------------------------------------------------------
# 2) Append the new cycle: rows 828-841, same 14-row pattern as always
#    (Name, URL, "Disponible", timestamp), each with its own hyperlink
#    on column B. Row 836 (MapStore) carries a "#/" fragment, which
#    Excel records as a hyperlink SubAddress ("location") instead of
#    folding it into the Target.
# ---------------------------------------------------------------------
$newStartRow = 828
$newTimestamp = 44232.91759111339

$entries = @(
    @("Odoo",               "https://www.dataintelligence-group.com/",                     ""),
    @("Blackbox",           "https://serviciodashboard.azurewebsites.net/",                 ""),
    @("PowerBI",            "https://powerbi.microsoft.com/es-es/",                         ""),
    @("Dropbox",            "https://www.dropbox.com/",                                     ""),
    @("Odoo",               "https://dataintelligence.store/",                              ""),
    @("GEE",                "https://app-data-i.users.earthengine.app/",                    ""),
    @("UtilidadesOdoo",     "https://odooutil.azurewebsites.net/",                           ""),
    @("Filtros Dashboard",  "https://filtradordashboard.azurewebsites.net/",                 ""),
    @("MapStore",           "https://ide.dataintelligence-group.com/mapstore/",              "/"),
    @("GeoServer",          "https://ide.dataintelligence-group.com/geoserver/web/?0",       ""),
    @("Tomcat",             "https://ide.dataintelligence-group.com/",                       ""),
    @("Shiny",              "https://rpubs.com/dataintelligence/",                           ""),
    @("Github",             "https://github.com/Sud-Austral/",                               ""),
    @("EZ Exporter",        "https://ezexporter.highviewapps.com/exports/export-profile/",   "")
)

# Carry over the A:D cell formatting (header-row style, hyperlink style,
# date style, ...) from the previous cycle before filling in the values,
# so the new rows pick up the same styles (s="2" on B, s="3" on D) as
# every other cycle in the sheet without inventing new style records.
$ws.Range("A814:D827").Copy() | Out-Null
$ws.Range("A828:D841").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $entries.Count; $i++) {
    $row = $newStartRow + $i
    $name = $entries[$i][0]
    $baseUrl = $entries[$i][1]
    $subAddress = $entries[$i][2]
    if ($subAddress -ne "") {
        $fullUrl = $baseUrl + "#" + $subAddress
    } else {
        $fullUrl = $baseUrl
    }

    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $fullUrl
    $ws.Cells.Item($row, 3).Value = "Disponible"
    $ws.Cells.Item($row, 4).Value = $newTimestamp

    if ($subAddress -ne "") {
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 2), $baseUrl, $subAddress) | Out-Null
    } else {
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 2), $baseUrl) | Out-Null
    }
}

# Hyperlinks.Add() re-applies formatting as it inserts the relationship;
# reassert the shared "Hyperlink" cell style on the whole new column-B
# block so it lines up with every other cycle in the sheet (s="2").
$ws.Range("B828:B841").Style = "Hyperlink"
